# Updated FRA model - 2025-08-10 22:21
# Sheet "VEDA_Sets-Proc" (second sheet / the active tab):
#  - B19: pattern string updated from "*bat*,-*EV*" to "EN*STG?hb*,-*EV*"
#  - H19/I19: add the T_Pos_AndOr / T_Neg_AndOr "And" / "Or" pair that the
#    other set rows (B3/B7/B17 group) already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

$ws.Range("B19").Value = "EN*STG?hb*,-*EV*"
$ws.Range("H19").Value = "And"
$ws.Range("I19").Value = "Or"
